# Update database and change read_price algorithm
# The sheet shows five consecutive twelve-month periods ending 1396/12..1400/12.
# This update rolls the window forward by one year: 1397/12..1401/12, shifting
# every numeric series one column to the left (E<-F, F<-G, G<-H, H<-I) and
# filling column I with the newly reported period's figures (some of the
# shifted-in H values were also revised).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column headers (row 8 and row 24): twelve-month period labels ---
$periods = @(
    "دوازده ماهه منتهی به 1397/12",
    "دوازده ماهه منتهی به 1398/12",
    "دوازده ماهه منتهی به 1399/12",
    "دوازده ماهه منتهی به 1400/12",
    "دوازده ماهه منتهی به 1401/12"
)
$cols = @(5, 6, 7, 8, 9)   # E, F, G, H, I

for ($i = 0; $i -lt 5; $i++) {
    $ws.Cells.Item(8, $cols[$i]).Value = $periods[$i]
    $ws.Cells.Item(24, $cols[$i]).Value = $periods[$i]
}

# --- Data rows: new E,F,G,H,I values per row ---
$data = @{
    10 = @(811858, 2065270, 2641204, 4656348, 7857597)
    11 = @(0, 0, 0, 0, 0)
    12 = @(326483, 390038, 600948, 1394825, 1267841)
    13 = @(546438, 860952, 1450615, 3146996, 5070000)
    14 = @(671100, 1352104, 0, 0, 0)
    15 = @(133522, 106049, 0, 0, 0)
    16 = @(198902, 234366, 266622, 326369, 789283)
    17 = @(2689601, 3912653, 6904322, 10917729, 16896066)
    18 = @(0, 0, 0, 0, 0)
    19 = @(2314633, 3016744, 6815422, 13885474, 21426500)
    20 = @(7692537, 11938176, 18679133, 34327741, 53307287)
    26 = @(2147, 1931, 2230, 2534, 1679)
    27 = @(10989, 10732, 9999, 9627, 10267)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    for ($i = 0; $i -lt 5; $i++) {
        $ws.Cells.Item($row, $cols[$i]).Value = $vals[$i]
    }
}
